# Convert Integer to the Sum of Two No-Zero Integers.
#
# 1. Row 64 col D currently holds the shared string "If-else/Switch-case";
#    append the extra note used by the updated pseudo-code sheet.
# 2. Append two brand-new problem rows (Id + Title only) at the bottom of
#    the table: row 65 ("Fruits Into Baskets 2") and row 66
#    ("Largest 3-Same-Digit Number in a String").
# 3. Leave the final selection on the last edited cell (B66), matching the
#    author's last on-screen action.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D64").Value = "If-else/Switch-case/ Use flag = (i % 3 == 0) + 2 * (i % 5 == 0);"

$ws.Range("A65").Value = 3477
$ws.Range("B65").Value = "Fruits Into Baskets 2"

$ws.Range("A66").Value = 2264
$ws.Range("B66").Value = "Largest 3-Same-Digit Number in a String"

$ws.Range("B66").Select() | Out-Null
